$wb = $excel.ActiveWorkbook

# The "Status" value changes from "Ready for handoff" to "In Translation"
# everywhere it appears: the Overview sheet's zh-cn/de-de status columns
# (E2/F2) and the per-language "Status" column (C2) on both the zh-cn and
# de-de detail sheets. Shrinking that text causes Excel to re-fit the
# width of the corresponding "Status" columns on each sheet.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E1").ColumnWidth = 12.5
$wsOverview.Range("F1").ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C1").ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C1").ColumnWidth = 12.5
